$wb = $excel.ActiveWorkbook

# Sheet "Retornos": remove the 4 most recent rows of data (rows 2:5),
# shifting all remaining data rows up by 4.
$ws1 = $wb.Worksheets.Item("Retornos")
$ws1.Range("A2:C5").EntireRow.Delete()

# Sheet "Beta": update the computed beta value to reflect the new date range.
$ws2 = $wb.Worksheets.Item("Beta")
$ws2.Range("B2").Value = -0.08241950177914951
